$d = $word.ActiveDocument

# 1) Strike through the first bullet item:
#    "Urediti da se unaprijed odabere 10 nasumicnih monstera u listi "
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Urediti da se unaprijed odabere 10 nasumicnih monstera u listi")) {
        $p.Range.Font.StrikeThrough = 1
        break
    }
}

# 2) Append two new bullet paragraphs after the last paragraph
#    (same ListParagraph style / numId=2 list as the rest of the list)
$lastP = $d.Paragraphs.Last
$lastP.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "Mozda bi trebalo dodat da user na pocetku moze izabrat svoj health i damage"

$p1 = $d.Paragraphs.Last
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = " "
